$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2299.102676639999
$ws.Range("E2").Value = 249384.345679256
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 138781.1577882601
$ws.Range("L2").Value = 449634.7497589202
$ws.Range("M2").Value = 101123.794939365
$ws.Range("N2").Value = 64159.90086799784
$ws.Range("O2").Value = 60521.71677636998

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 3009.804799738145
$ws.Range("B2").Value = 27164.01933074473
$ws.Range("E2").Value = 148029.3904973308
$ws.Range("I2").Value = 167745.2736830803
$ws.Range("L2").Value = 40177.05802310962
$ws.Range("M2").Value = 55411.58404378576
$ws.Range("N2").Value = 17452.07979808242
$ws.Range("O2").Value = 10096.07155031569

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 18346.45888886983
$ws.Range("B2").Value = 15997.19565778882
$ws.Range("E2").Value = 107448.9872383075
$ws.Range("I2").Value = 151988.1037676829
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 55328.40432004048
$ws.Range("N2").Value = 39263.10978901467
$ws.Range("O2").Value = 46419.8816816441
